$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 299, pushing the existing rows 299-310
# (which contained the historical data) down to become rows 301-312.
$ws.Rows.Item(299).Insert()
$ws.Rows.Item(299).Insert()

# Populate the first new row (299) with the latest "1a amarillo" record.
$ws.Cells.Item(299, 1).Value = 4
$ws.Cells.Item(299, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(299, 3).Value = "Los Lagos"
$ws.Cells.Item(299, 4).Value = 44509
$ws.Cells.Item(299, 5).Value = 10
$ws.Cells.Item(299, 6).Value = "Fruta"
$ws.Cells.Item(299, 7).Value = 100102
$ws.Cells.Item(299, 8).Value = "Cítricos"
$ws.Cells.Item(299, 9).Value = 100102003
$ws.Cells.Item(299, 10).Value = "Limón"
$ws.Cells.Item(299, 11).Value = "Sin especificar"
$ws.Cells.Item(299, 12).Value = "1a amarillo"
$ws.Cells.Item(299, 13).Value = 1200
$ws.Cells.Item(299, 14).Value = 11000
$ws.Cells.Item(299, 15).Value = 12000
$ws.Cells.Item(299, 16).Value = 11500
$ws.Cells.Item(299, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(299, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(299, 19).Value = 639
$ws.Cells.Item(299, 20).Value = 18

# Populate the second new row (300) with the latest "2a amarillo" record.
$ws.Cells.Item(300, 1).Value = 4
$ws.Cells.Item(300, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(300, 3).Value = "Los Lagos"
$ws.Cells.Item(300, 4).Value = 44509
$ws.Cells.Item(300, 5).Value = 10
$ws.Cells.Item(300, 6).Value = "Fruta"
$ws.Cells.Item(300, 7).Value = 100102
$ws.Cells.Item(300, 8).Value = "Cítricos"
$ws.Cells.Item(300, 9).Value = 100102003
$ws.Cells.Item(300, 10).Value = "Limón"
$ws.Cells.Item(300, 11).Value = "Sin especificar"
$ws.Cells.Item(300, 12).Value = "2a amarillo"
$ws.Cells.Item(300, 13).Value = 400
$ws.Cells.Item(300, 14).Value = 9000
$ws.Cells.Item(300, 15).Value = 9000
$ws.Cells.Item(300, 16).Value = 9000
$ws.Cells.Item(300, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(300, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(300, 19).Value = 500
$ws.Cells.Item(300, 20).Value = 18
